# Generate Report for Handback
#
# This applies the "handback" update to the localization-status workbook:
#   - Overview!C2 status text updated (via the shared "Ready for handoff"
#     string, which is also used on the zh-cn/de-de sheets).
#   - zh-cn / de-de sheets: the "Latest Target File" (I) and
#     "Latest Handback File" (J) columns get filled in for rows 2 & 3, with a
#     hyperlink added on the "Latest Target File" cell (matching the existing
#     hyperlink on column A).
#   - zh-cn's "Latest Handback DateTime" (K) timestamp is updated; de-de's is
#     set to a new (different) timestamp.
#   - Column widths are widened on both the Overview and per-language sheets
#     to better fit the newly-populated content.

$wb = $excel.ActiveWorkbook

$mdDisplay = "bd32cb5f-a5de-46e6-b32d-cc1b6cc08662.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c56efae525b0f08417e331f2570346b140e490f2/e2e/bd32cb5f-a5de-46e6-b32d-cc1b6cc08662.md"

# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared across Overview!E2:F3 and the Status column (C)
#    on both language sheets, so updating every cell that currently holds it
#    keeps the shared string table in sync (and matches what Excel would do
#    editing any one occurrence).
#    NOTE: cells are located first (read-only pass) and updated in a second
#    pass - mutating values while enumerating a live Cells collection is
#    unreliable. Also, string literals are kept on the left of -eq so
#    PowerShell doesn't coerce Boolean-valued cell text (e.g. "True") to
#    $true and match every non-empty comparison.
$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $matchAddrs = @()
    foreach ($cell in $used.Cells) {
        if ($statusOld -eq $cell.Text) {
            $matchAddrs += $cell.Address()
        }
    }
    foreach ($addr in $matchAddrs) {
        $ws.Range($addr).Value = $statusNew
    }
}

# 2. Per-language sheets: fill in Latest Target File (I) / Latest Handback
#    File (J) / Latest Handback DateTime (K) for rows 2 and 3, and add the
#    matching hyperlink on column I (same display text/target as column A's
#    existing hyperlink).
function Set-HandbackRow($ws, $row, $handbackFile, $handbackDate) {
    $ws.Range("I$row").Value = $mdDisplay
    $ws.Hyperlinks.Add($ws.Range("I$row"), $mdUrl, "", "", $mdDisplay)
    $ws.Range("J$row").Value = $handbackFile
    $ws.Range("K$row").Value = $handbackDate
}

$wsZh = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZh 2 "bd32cb5f-a5de-46e6-b32d-cc1b6cc08662.be67c59b958f36f60f42e30976160bd87187b347.zh-cn.xlf" "2016-08-25 17:05:42"
Set-HandbackRow $wsZh 3 "bd32cb5f-a5de-46e6-b32d-cc1b6cc08662.be67c59b958f36f60f42e30976160bd87187b347.zh-cn.xlf" "2016-08-25 17:05:42"

$wsDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDe 2 "bd32cb5f-a5de-46e6-b32d-cc1b6cc08662.be67c59b958f36f60f42e30976160bd87187b347.de-de.xlf" "2016-08-25 17:05:49"
Set-HandbackRow $wsDe 3 "bd32cb5f-a5de-46e6-b32d-cc1b6cc08662.be67c59b958f36f60f42e30976160bd87187b347.de-de.xlf" "2016-08-25 17:05:49"

# 3. Widen columns to fit the newly-populated / longer content.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

foreach ($ws in @($wsZh, $wsDe)) {
    $ws.Range("C1").ColumnWidth = 29.9777047293527
    $ws.Range("I1").ColumnWidth = 40
    $ws.Range("J1").ColumnWidth = 40
}
